$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on price cells that would otherwise be auto-coerced
# to numbers by Excel (values that look like plain decimals).
$textCells = @("D5","D6","D10","D16","D19","D20","D21","D22","D23","D24","D25","D27","D29","D32","D33","D34","D36","D38","D39","D40","D41","D42","D43","D44","D45","D47","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '60.230.39'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.617.21'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '522.11'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').Value = '149.26'
$ws.Range('E6').Value = '  -2.80%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -4.44%  '
$ws.Range('D9').Value = '2.619.86'
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('D10').Value = '6.28'
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').Value = '3.070.07'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '60.215.77'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '21.22'
$ws.Range('E16').Value = '  -2.48%  '
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').Value = '2.614.80'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '4.64'
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '341.72'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').Value = '10.41'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').Value = '6.11'
$ws.Range('E22').Value = '  -1.83%  '
$ws.Range('D23').Value = '0.995'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '60.64'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '0.420'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('D27').Value = '0.162'
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').Value = '0.0₃0809'
$ws.Range('E28').Value = '  -3.88%  '
$ws.Range('D29').Value = '7.08'
$ws.Range('E29').Value = '  -3.51%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('D32').Value = '1.59'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '18.97'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').Value = '149.31'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').Value = '0.919'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('E37').Value = '  -5.02%  '
$ws.Range('D38').Value = '0.864'
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('D39').Value = '36.50'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').Value = '1.44'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('D41').Value = '3.63'
$ws.Range('E41').Value = '  -4.25%  '
$ws.Range('D42').Value = '290.25'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.624'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.100'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -2.21%  '
$ws.Range('D47').Value = '19.48'
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').Value = '4.70'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('D51').Value = '1.955.04'
$ws.Range('E51').Value = '  -0.91%  '
